$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("uf-tx-sucesso")

# Shift the "MA"/"MT" labels: A3 becomes "MA", A4 becomes "MT", A5 becomes empty.
$ws.Range("A3").Value = "MA"
$ws.Range("A4").Value = "MT"
$ws.Range("A5").ClearContents()
